$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cut-off date (Fecha_Corte). C3:C11 are formulas chained off C2
# (=C2, =C3, ...) so updating C2 propagates automatically.
$ws.Range("C2").Value = 46071

# Updated policy counts (Polizas_Totales) and commissions (Comisones) per
# advisor row, reflecting the new Feb 18 2026 cut-off data.
$ws.Range("H3").Value = 45727.57

$ws.Range("H5").Value = 69639.95

$ws.Range("G6").Value = 10.5

$ws.Range("H8").Value = 31409.37

$ws.Range("H9").Value = 85425.04

$ws.Range("G10").Value = 9.5
$ws.Range("H10").Value = 49591.66

$ws.Range("G11").Value = 8.5
$ws.Range("H11").Value = 38447.01

# Move the active selection on the sheet (view-state only).
$ws.Range("F18").Select()
